$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix a typo in the existing mapping table (row 8, column B) ---
$ws.Cells.Item(8, 2).Value = "Strategic Accounts Marketing"

# --- Bold the header row (A1:B1) ---
$ws.Range("A1:B1").Font.Bold = $true

# --- Append the new "Lead Source" mapping rows (rows 10-24) ---
$newRows = @(
    @("Chat", "Chat"),
    @("Content Syndication", "Content Syndication"),
    @("Email", "Email"),
    @("enewal", "enewal"),
    @("Events", "Events"),
    @("Inc. 5", "Inc. 6"),
    @("InsideView", "InsideView"),
    @("Jigsaw", "Jigsaw"),
    @("LinkedIn - Outbound", "LinkedIn - Outbound"),
    @("Online Events", "Online Events"),
    @("Outbound", "Outbound"),
    @("Physical Event", "Physical Event"),
    @("Physical Events", "Physical Event"),
    @("Social", "Social"),
    @("Website", "Website")
)

$startRow = 10
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $pair = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $ws.Cells.Item($r, 2).Value = $pair[1]
}

# --- Update the selected cell to match the saved view state ---
$ws.Range("B12").Select() | Out-Null
